# Append new match rows (12-21) to the "Kieron Pollard " worksheet.
# Data mirrors the existing rows (2-11), re-scraped in a different order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "batsman" is "Kieron Pollard" followed by a trailing U+00A0 (non-breaking
# space), matching the source data exactly.
$nbsp = [char]0x00A0
$batsman = "Kieron Pollard$nbsp"

$rows = @(
    @(" Abu Dhabi", " September 19 2020", "Super Kings won by 5 wickets (with 4 balls remaining)", "Mumbai Indians", "Chennai Super Kings", $batsman, "18", "14", "1", "1", "128.57"),
    @(" Abu Dhabi", " September 23 2020", "Mumbai won by 49 runs", "Mumbai Indians", "Kolkata Knight Riders", $batsman, "13", "7", "1", "0", "185.71"),
    @(" Abu Dhabi", " October 11 2020", "Mumbai won by 5 wickets (with 2 balls remaining)", "Mumbai Indians", "Delhi Capitals", $batsman, "11", "14", "1", "0", "78.57"),
    @(" Dubai (DSC)", " November 05 2020", "Mumbai won by 57 runs", "Mumbai Indians", "Delhi Capitals", $batsman, "0", "2", "0", "0", "0.00"),
    @(" Abu Dhabi", " October 01 2020", "Mumbai won by 48 runs", "Mumbai Indians", "Kings XI Punjab", $batsman, "47", "20", "3", "4", "235.00"),
    @(" Sharjah", " October 04 2020", "Mumbai won by 34 runs", "Mumbai Indians", "Sunrisers Hyderabad", $batsman, "25", "13", "0", "3", "192.30"),
    @(" Sharjah", " November 03 2020", "Sunrisers won by 10 wickets (with 17 balls remaining)", "Mumbai Indians", "Sunrisers Hyderabad", $batsman, "41", "25", "2", "4", "164.00"),
    @(" Dubai (DSC)", " November 10 2020", "Mumbai won by 5 wickets (with 8 balls remaining)", "Mumbai Indians", "Delhi Capitals", $batsman, "9", "4", "2", "0", "225.00"),
    @(" Dubai (DSC)", " September 28 2020", "Match tied (RCB won the one-over eliminator)", "Mumbai Indians", "Royal Challengers Bangalore", $batsman, "60", "24", "3", "5", "250.00"),
    @(" Dubai (DSC)", " October 18 2020", "Match tied (Kings XI won the one-over eliminator)", "Mumbai Indians", "Kings XI Punjab", $batsman, "34", "12", "1", "4", "283.33")
)

# Columns G,H,I,J,K (index 7-11) hold numeric-looking text (totalRuns, totalBalls,
# total4s, total6s, sr) that must stay text, like the rest of the sheet - so they
# are entered with a leading apostrophe to force text storage instead of Excel's
# automatic number conversion.
$textCols = @(7, 8, 9, 10, 11)

$startRow = 12
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 1; $c -le $row.Count; $c++) {
        $value = $row[$c - 1]
        if ($textCols -contains $c) {
            $ws.Cells.Item($r, $c).Value = "'" + $value
        } else {
            $ws.Cells.Item($r, $c).Value = $value
        }
    }
}
